$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 126, pushing the existing rows
# 126-180 down to 128-182 (dimension grows from R180 to R182).
$ws.Rows.Item(126).Insert()
$ws.Rows.Item(126).Insert()

# New row 126: weekly Berenjena record for 2021-09-27 (serial 44466)
$ws.Range("A126").Value = 10
$ws.Range("B126").Value = "Vega Modelo de Temuco"
$ws.Range("C126").Value = "La Araucanía"
$ws.Range("D126").Value = 44466
$ws.Range("E126").Value = 9
$ws.Range("F126").Value = 100112001
$ws.Range("G126").Value = "Berenjena"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 150
$ws.Range("K126").Value = 10000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = 11333
$ws.Range("N126").Value = "$/caja 60 unidades"
$ws.Range("O126").Value = "Región de Arica y Parinacota"
$ws.Range("P126").Value = 189
$ws.Range("Q126").Value = 60
$ws.Range("R126").Value = "Hortaliza"

# New row 127: weekly Berenjena record for 2021-09-27 (serial 44466)
$ws.Range("A127").Value = 10
$ws.Range("B127").Value = "Vega Modelo de Temuco"
$ws.Range("C127").Value = "La Araucanía"
$ws.Range("D127").Value = 44466
$ws.Range("E127").Value = 9
$ws.Range("F127").Value = 100112001
$ws.Range("G127").Value = "Berenjena"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Segunda"
$ws.Range("J127").Value = 20
$ws.Range("K127").Value = 8000
$ws.Range("L127").Value = 8000
$ws.Range("M127").Value = 8000
$ws.Range("N127").Value = "$/caja 90 unidades"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 89
$ws.Range("Q127").Value = 90
$ws.Range("R127").Value = "Hortaliza"
